$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @('D2', '242.76'),
  @('G2', '19'),
  @('G3', '19'),
  @('D4', '5.388'),
  @('G4', '19'),
  @('D5', '0.05939'),
  @('G5', '19'),
  @('D6', '3.400'),
  @('G6', '19'),
  @('B7', 'KuCoinToken'),
  @('C7', 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'),
  @('D7', '6.450'),
  @('E7', '6KuCoinTokenKCS'),
  @('G7', '19'),
  @('B8', 'MXToken'),
  @('C8', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
  @('D8', '0.8068'),
  @('E8', '7MXTokenMX'),
  @('G8', '19'),
  @('B9', 'FTXToken'),
  @('C9', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
  @('D9', '0.9031'),
  @('E9', '8FTXTokenFTT'),
  @('G9', '19'),
  @('B10', 'WazirX'),
  @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
  @('D10', '0.1418'),
  @('E10', '9WazirXWRX'),
  @('G10', '19'),
  @('B11', 'MandalaExchangeToken'),
  @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
  @('D11', '0.07443'),
  @('E11', '10MandalaExchangeTokenMDX'),
  @('G11', '19'),
  @('B12', 'LiechtensteinCryptoassetsExchange'),
  @('C12', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
  @('D12', '0.03261'),
  @('E12', '11LiechtensteinCryptoassetsExchangeLCX'),
  @('G12', '19'),
  @('B13', 'BitrueCoin'),
  @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
  @('D13', '0.03044'),
  @('E13', '12BitrueCoinBTR'),
  @('G13', '19'),
  @('B14', 'BitMartToken'),
  @('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
  @('D14', '0.09332'),
  @('E14', '13BitMartTokenBMX'),
  @('G14', '19'),
  @('B15', 'MCDex'),
  @('C15', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'),
  @('D15', '3.950'),
  @('E15', '14MCDexMCB'),
  @('G15', '19'),
  @('B16', 'BitForexToken'),
  @('C16', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
  @('D16', '0.001578'),
  @('E16', '15BitForexTokenBF'),
  @('G16', '19'),
  @('B17', 'CoinExToken'),
  @('C17', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
  @('D17', '0.04793'),
  @('E17', '16CoinExTokenCET'),
  @('G17', '19'),
  @('B18', 'One'),
  @('C18', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
  @('D18', '0.0005948'),
  @('E18', '17OneONE'),
  @('G18', '19'),
  @('B19', 'TigerCash'),
  @('C19', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
  @('D19', '0.006311'),
  @('E19', '18TigerCashTCH'),
  @('G19', '19'),
  @('B20', 'UpBots'),
  @('C20', 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'),
  @('D20', '0.007493'),
  @('E20', '19UpBotsUBXTBestin24h'),
  @('G20', '19'),
  @('B21', 'HotbitToken'),
  @('C21', 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'),
  @('D21', '0.004408'),
  @('E21', '20HotbitTokenHTB'),
  @('G21', '19'),
  @('B22', 'BitKan'),
  @('C22', 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'),
  @('D22', '0.0009880'),
  @('E22', '21BitKanKAN'),
  @('G22', '19'),
  @('B23', 'NitroEx'),
  @('C23', 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'),
  @('D23', '0.00007811'),
  @('E23', '22NitroExNTX'),
  @('G23', '19'),
  @('B24', 'LEO'),
  @('C24', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
  @('D24', '3.613'),
  @('E24', '23LEOLEO'),
  @('G24', '19'),
  @('G25', '19'),
  @('G26', '19'),
  @('D27', '0.1322'),
  @('G27', '19'),
  @('G28', '19'),
  @('G29', '19'),
  @('G30', '19'),
  @('G31', '19'),
  @('G32', '19'),
  @('G33', '19'),
  @('G34', '19'),
  @('G35', '19'),
  @('G36', '19'),
  @('G37', '19'),
  @('G38', '19'),
  @('G39', '19'),
  @('D40', '0.03876'),
  @('G40', '19'),
  @('D41', '0.006229'),
  @('G41', '19'),
  @('D42', '0.1068'),
  @('G42', '19'),
  @('D43', '0.002614'),
  @('G43', '19'),
  @('D44', '0.007273'),
  @('G44', '19'),
  @('D45', '0.00005198'),
  @('G45', '19'),
  @('D46', '0.00000000751'),
  @('G46', '19'),
  @('D47', '0.0005807'),
  @('G47', '19'),
  @('D48', '0.9121'),
  @('G48', '19'),
  @('D49', '0.002263'),
  @('G49', '19'),
  @('D50', '0.00002103'),
  @('G50', '19'),
  @('D51', '0.0002003'),
  @('G51', '19')
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newVal = $pair[1]
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newVal
    $c.ClearFormats()
}
